$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '58.310.99'
$ws.Range('E2').Value = '  -1.04%  '

$ws.Range('D3').Value = '2.476.98'
$ws.Range('E3').Value = '  -0.99%  '

$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '0.999'
$ws.Range('E4').Value = '  -0.13%  '

$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '520.21'
$ws.Range('E5').Value = '  -2.25%  '

$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '134.05'
$ws.Range('E6').Value = '  -0.89%  '

$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.999'
$ws.Range('E7').Value = '  -0.12%  '

$ws.Range('E8').Value = '  -1.38%  '

$ws.Range('D9').Value = '2.491.51'
$ws.Range('E9').Value = '  -0.54%  '

$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.0985'
$ws.Range('E10').Value = '  -2.79%  '

$ws.Range('E11').Value = '  -0.91%  '

$ws.Range('E12').Value = '  -1.75%  '

$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '0.337'
$ws.Range('E13').Value = '  -2.19%  '

$ws.Range('D14').Value = '2.920.80'
$ws.Range('E14').Value = '  -0.87%  '

$ws.Range('D15').Value = '58.189.65'
$ws.Range('E15').Value = '  -1.12%  '

$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '22.04'
$ws.Range('E16').Value = '  -3.00%  '

$ws.Range('E17').Value = '  -1.94%  '

$ws.Range('D18').Value = '2.493.05'
$ws.Range('E18').Value = '  -0.64%  '

$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '10.64'
$ws.Range('E19').Value = '  -3.40%  '

$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '321.23'
$ws.Range('E20').Value = '  -0.67%  '

$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '4.17'
$ws.Range('E21').Value = '  -1.44%  '

$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '0.999'
$ws.Range('E22').Value = '  -0.03%  '

$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '5.75'
$ws.Range('E23').Value = '  -2.85%  '

$ws.Range('E24').Value = '  -0.16%  '

$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '0.411'
$ws.Range('E25').Value = '  -1.89%  '

$ws.Range('B26').Value = 'Binance-PegBSC-USD'
$ws.Range('C26').Value = 'https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd'
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '0.997'
$ws.Range('E26').Value = '  -0.23%  '

$ws.Range('B27').Value = 'Kaspa'
$ws.Range('C27').Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '0.161'
$ws.Range('E27').Value = '  -1.26%  '

$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '7.37'
$ws.Range('E28').Value = '  -1.69%  '

$ws.Range('D29').Value = '0.0₃0749'
$ws.Range('E29').Value = '  -1.53%  '

$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '169.62'
$ws.Range('E30').Value = '  +0.00%  '

$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '6.33'
$ws.Range('E31').Value = '  -1.80%  '

$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '1.69'
$ws.Range('E32').Value = '  -2.61%  '

$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '1.18'
$ws.Range('E33').Value = '  +4.49%  '

$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '0.998'
$ws.Range('E34').Value = '  -0.04%  '

$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '0.998'
$ws.Range('E35').Value = '  -0.04%  '

$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '18.05'
$ws.Range('E36').Value = '  -1.45%  '

$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '1.31'
$ws.Range('E37').Value = '  -3.31%  '

$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '4.01'
$ws.Range('E38').Value = '  -0.42%  '

$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '36.64'
$ws.Range('E39').Value = '  -0.47%  '

$ws.Range('E40').Value = '  -3.20%  '

$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '0.799'
$ws.Range('E41').Value = '  +0.50%  '

$ws.Range('B42').Value = 'RenderToken'
$ws.Range('C42').Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '5.12'
$ws.Range('E42').Value = '  +3.06%  '

$ws.Range('B43').Value = 'Bittensor'
$ws.Range('C43').Value = 'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao'
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '275.04'
$ws.Range('E43').Value = '  -1.79%  '

$ws.Range('B44').Value = 'Filecoin'
$ws.Range('C44').Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '3.44'
$ws.Range('E44').Value = '  -3.62%  '

$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '0.598'
$ws.Range('E45').Value = '  -0.60%  '

$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '123.95'
$ws.Range('E46').Value = '  -3.88%  '

$ws.Range('E47').Value = '  -1.43%  '

$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '0.0490'
$ws.Range('E48').Value = '  -1.83%  '

$ws.Range('E49').Value = '  -1.80%  '

$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '17.07'
$ws.Range('E50').Value = '  -0.84%  '

$ws.Range('D51').Value = '1.736.51'
$ws.Range('E51').Value = '  -0.58%  '
